# Wave 5_React_Tracker_(Pandian).xlsx update
# - Header row (row 3): "Static Project Layout" moves earlier, "Labels and
#   Inbox Mails" header is split into separate "Labels" / "Inbox Mails"
#   headers, and a new trailing "React Router" column (L) is added.
# - Data row (row 4): A4 now holds an actual date, columns B-F are marked
#   "done", G4 is marked "in progress" and H4 stays blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3: header labels ----
$ws.Range("A3").Value = "Date"
$ws.Range("B3").Value = "React Basics (Architecture, Concepts, Lifecycle)"
$ws.Range("C3").Value = "Properties Flow"
$ws.Range("D3").Value = "State Change"
$ws.Range("E3").Value = "Static Project Layout"
$ws.Range("G3").Value = "Inbox Mails"
$ws.Range("F3").Value = "Labels"
$ws.Range("H3").Value = "Label Functionality (SENT/IMPORTANT/DRAFT)"
$ws.Range("I3").Value = "Compose Mail"
$ws.Range("J3").Value = "View Mail"
$ws.Range("K3").Value = "Reply"
$ws.Range("L3").Value = "React Router"
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Row 4: data values ----
$ws.Range("A4").Value = "23/08/2016"
$ws.Range("B4").Value = "done"
$ws.Range("C4").Value = "done"
$ws.Range("D4").Value = "done"
$ws.Range("E4").Value = "done"
$ws.Range("F4").Value = "done"
$ws.Range("G4").Value = "in progress"

# ---- Selection ----
$ws.Range("G4").Select() | Out-Null
